# Weekly refresh of fruit/vegetable price rows: the source rows 4-26
# (excluding 6-7, which are untouched) get their Fecha/Volumen/Precio
# columns re-shuffled to reflect the latest weekly pull. Apply the new
# per-row values directly (destination row -> new D/J/K/L/M/P values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @{}
$newData[26] = @{ D = 44474; J = 40;  K = 13000; L = 14000; M = 13500; P = 1038 }
$newData[17] = @{ D = 44362; J = 40;  K = 15000; L = 16000; M = 15500; P = 1192 }
$newData[16] = @{ D = 44159; J = 60;  K = 30000; L = 32000; M = 31000; P = 2385 }
$newData[12] = @{ D = 44435; J = 100; K = 13000; L = 14000; M = 13500; P = 1038 }
$newData[14] = @{ D = 44433; J = 100; K = 13000; L = 14000; M = 13500; P = 1038 }
$newData[13] = @{ D = 44350; J = 40;  K = 23000; L = 25000; M = 24000; P = 1846 }
$newData[25] = @{ D = 44453; J = 50;  K = 14000; L = 15000; M = 14600; P = 1123 }
$newData[23] = @{ D = 44523; J = 40;  K = 15000; L = 16000; M = 15500; P = 1192 }
$newData[19] = @{ D = 44510; J = 40;  K = 15000; L = 16000; M = 15500; P = 1192 }
$newData[11] = @{ D = 44308; J = 50;  K = 26000; L = 27000; M = 26400; P = 2031 }
$newData[5]  = @{ D = 44320; J = 50;  K = 26000; L = 28000; M = 26800; P = 2062 }
$newData[24] = @{ D = 44316; J = 50;  K = 27000; L = 28000; M = 27400; P = 2108 }
$newData[8]  = @{ D = 44467; J = 100; K = 13000; L = 14000; M = 13500; P = 1038 }
$newData[20] = @{ D = 44313; J = 50;  K = 25000; L = 26000; M = 25600; P = 1969 }
$newData[21] = @{ D = 44334; J = 50;  K = 26000; L = 28000; M = 27200; P = 2092 }
$newData[10] = @{ D = 44509; J = 100; K = 15000; L = 16000; M = 15500; P = 1192 }
$newData[9]  = @{ D = 44488; J = 40;  K = 16000; L = 17000; M = 16500; P = 1269 }
$newData[22] = @{ D = 44264; J = 40;  K = 30000; L = 32000; M = 31000; P = 2385 }
$newData[15] = @{ D = 44327; J = 50;  K = 24000; L = 25000; M = 24400; P = 1877 }
$newData[4]  = @{ D = 44503; J = 35;  K = 15000; L = 16000; M = 15429; P = 1187 }
$newData[18] = @{ D = 44462; J = 60;  K = 14000; L = 15000; M = 14500; P = 1115 }

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
